$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 38.04655766666666
$ws.Range("H2").Value = 114.139673
$ws.Range("I2").Value = 0.8090698722086991
$ws.Range("J2").Value = 0.8090698722086992
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 17.13024733333333
$ws.Range("N2").Value = 51.390742
$ws.Range("O2").Value = 0.2959211466465044
$ws.Range("P2").Value = 0.2959211466465043
$ws.Range("Q2").Value = 651.7469430119295
$ws.Range("R2").Value = 5865.722487107366
$ws.Range("S2").Value = 0.239420884301139
$ws.Range("T2").Value = 0.239420884301139

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 38.04655766666666
$ws.Range("H3").Value = 114.139673
$ws.Range("I3").Value = 0.8090698722086991
$ws.Range("J3").Value = 0.8090698722086992
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 27.61090666666666
$ws.Range("N3").Value = 82.83272
$ws.Range("O3").Value = 0.4769721651858779
$ws.Range("P3").Value = 0.4769721651858778
$ws.Range("Q3").Value = 1050.499952722284
$ws.Range("R3").Value = 9454.49957450056
$ws.Range("S3").Value = 0.3859038087340447
$ws.Range("T3").Value = 0.3859038087340447

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 38.04655766666666
$ws.Range("H4").Value = 114.139673
$ws.Range("I4").Value = 0.8090698722086991
$ws.Range("J4").Value = 0.8090698722086992
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 4.423514
$ws.Range("N4").Value = 13.270542
$ws.Range("O4").Value = 0.07641520344774541
$ws.Range("P4").Value = 0.0764152034477454
$ws.Range("Q4").Value = 168.2994804903073
$ws.Range("R4").Value = 1514.695324412766
$ws.Range("S4").Value = 0.06182523888826912
$ws.Range("T4").Value = 0.06182523888826912

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 38.04655766666666
$ws.Range("H5").Value = 114.139673
$ws.Range("I5").Value = 0.8090698722086991
$ws.Range("J5").Value = 0.8090698722086992
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 8.723210333333332
$ws.Range("N5").Value = 26.169631
$ws.Range("O5").Value = 0.1506914847198724
$ws.Range("P5").Value = 0.1506914847198724
$ws.Range("Q5").Value = 331.8881249856292
$ws.Range("R5").Value = 2986.993124870663
$ws.Range("S5").Value = 0.1219199402852463
$ws.Range("T5").Value = 0.1219199402852463

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1.617245333333334
$ws.Range("H6").Value = 4.851736000000001
$ws.Range("I6").Value = 0.03439113957782537
$ws.Range("J6").Value = 0.03439113957782537
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 17.13024733333333
$ws.Range("N6").Value = 51.390742
$ws.Range("O6").Value = 0.2959211466465044
$ws.Range("P6").Value = 0.2959211466465043
$ws.Range("Q6").Value = 27.70381255867911
$ws.Range("R6").Value = 249.334313028112
$ws.Range("S6").Value = 0.01017706545835006
$ws.Range("T6").Value = 0.01017706545835006

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1.617245333333334
$ws.Range("H7").Value = 4.851736000000001
$ws.Range("I7").Value = 0.03439113957782537
$ws.Range("J7").Value = 0.03439113957782537
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 27.61090666666666
$ws.Range("N7").Value = 82.83272
$ws.Range("O7").Value = 0.4769721651858779
$ws.Range("P7").Value = 0.4769721651858778
$ws.Range("Q7").Value = 44.65360995576889
$ws.Range("R7").Value = 401.88248960192
$ws.Range("S7").Value = 0.01640361630764511
$ws.Range("T7").Value = 0.0164036163076451

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 1.617245333333334
$ws.Range("H8").Value = 4.851736000000001
$ws.Range("I8").Value = 0.03439113957782537
$ws.Range("J8").Value = 0.03439113957782537
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 4.423514
$ws.Range("N8").Value = 13.270542
$ws.Range("O8").Value = 0.07641520344774541
$ws.Range("P8").Value = 0.0764152034477454
$ws.Range("Q8").Value = 7.153907373434668
$ws.Range("R8").Value = 64.38516636091201
$ws.Range("S8").Value = 0.002628005927639335
$ws.Range("T8").Value = 0.002628005927639334

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 1.617245333333334
$ws.Range("H9").Value = 4.851736000000001
$ws.Range("I9").Value = 0.03439113957782537
$ws.Range("J9").Value = 0.03439113957782537
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 8.723210333333332
$ws.Range("N9").Value = 26.169631
$ws.Range("O9").Value = 0.1506914847198724
$ws.Range("P9").Value = 0.1506914847198724
$ws.Range("Q9").Value = 14.10757120326845
$ws.Range("R9").Value = 126.968140829416
$ws.Range("S9").Value = 0.005182451884190872
$ws.Range("T9").Value = 0.005182451884190871

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 7.361255
$ws.Range("H10").Value = 22.083765
$ws.Range("I10").Value = 0.1565389882134754
$ws.Range("J10").Value = 0.1565389882134754
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 17.13024733333333
$ws.Range("N10").Value = 51.390742
$ws.Range("O10").Value = 0.2959211466465044
$ws.Range("P10").Value = 0.2959211466465043
$ws.Range("Q10").Value = 126.1001188337367
$ws.Range("R10").Value = 1134.90106950363
$ws.Range("S10").Value = 0.04632319688701529
$ws.Range("T10").Value = 0.04632319688701528

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 7.361255
$ws.Range("H11").Value = 22.083765
$ws.Range("I11").Value = 0.1565389882134754
$ws.Range("J11").Value = 0.1565389882134754
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 27.61090666666666
$ws.Range("N11").Value = 82.83272
$ws.Range("O11").Value = 0.4769721651858779
$ws.Range("P11").Value = 0.4769721651858778
$ws.Range("Q11").Value = 203.2509247545333
$ws.Range("R11").Value = 1829.2583227908
$ws.Range("S11").Value = 0.074664740144188
$ws.Range("T11").Value = 0.07466474014418799

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 7.361255
$ws.Range("H12").Value = 22.083765
$ws.Range("I12").Value = 0.1565389882134754
$ws.Range("J12").Value = 0.1565389882134754
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 4.423514
$ws.Range("N12").Value = 13.270542
$ws.Range("O12").Value = 0.07641520344774541
$ws.Range("P12").Value = 0.0764152034477454
$ws.Range("Q12").Value = 32.56261455007
$ws.Range("R12").Value = 293.06353095063
$ws.Range("S12").Value = 0.01196195863183695
$ws.Range("T12").Value = 0.01196195863183695

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 7.361255
$ws.Range("H13").Value = 22.083765
$ws.Range("I13").Value = 0.1565389882134754
$ws.Range("J13").Value = 0.1565389882134754
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 8.723210333333332
$ws.Range("N13").Value = 26.169631
$ws.Range("O13").Value = 0.1506914847198724
$ws.Range("P13").Value = 0.1506914847198724
$ws.Range("Q13").Value = 64.21377568230166
$ws.Range("R13").Value = 577.923981140715
$ws.Range("S13").Value = 0.02358909255043522
$ws.Range("T13").Value = 0.02358909255043522
